# Apply the "Office Theme" design (replacing the former "Integral" / "Red Violet"
# theme) to the presentation's slide master theme (ppt/theme/theme1.xml),
# mirroring what happens in PowerPoint when a different built-in theme is
# selected from the Design tab: the twelve theme colour-scheme slots are
# replaced with the standard Office palette, and the theme/colour-scheme
# names are updated to match.

$p = $ppt.ActivePresentation
$m = $p.SlideMaster
$t = $m.Theme

# Rename the theme and its colour scheme to match the new "Office Theme".
$t.Name = "Office Theme"
$tcs = $t.ThemeColorScheme
$tcs.Name = "Office"

# Standard Office theme colours (dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink)
# expressed as COM RGB() values (0x00BBGGRR).
$tcs.Colors(1).RGB = 0          # dk1      000000
$tcs.Colors(2).RGB = 16777215   # lt1      FFFFFF
$tcs.Colors(3).RGB = 6968388    # dk2      44546A
$tcs.Colors(4).RGB = 15132391   # lt2      E7E6E6
$tcs.Colors(5).RGB = 13998939   # accent1  5B9BD5
$tcs.Colors(6).RGB = 3243501    # accent2  ED7D31
$tcs.Colors(7).RGB = 10855845   # accent3  A5A5A5
$tcs.Colors(8).RGB = 49407      # accent4  FFC000
$tcs.Colors(9).RGB = 12874308   # accent5  4472C4
$tcs.Colors(10).RGB = 4697456   # accent6  70AD47
$tcs.Colors(11).RGB = 12673797  # hlink    0563C1
$tcs.Colors(12).RGB = 7491477   # folHlink 954F72
